$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("J2").Value = 1.05
$ws.Range("L2").Value = 1.37
$ws.Range("P2").Value = 1.5

# Row 3 updates
$ws.Range("J3").Value = 1.08
$ws.Range("L3").Value = 1.36
$ws.Range("P3").Value = 1.44
$ws.Range("Q3").Value = 2.63
$ws.Range("R3").Value = 1.8
$ws.Range("S3").Value = 1.8

# Row 4 updates - fill in odds that were previously empty
$ws.Range("G4").Value = 2.15
$ws.Range("H4").Value = 2.95
$ws.Range("I4").Value = 3.5
$ws.Range("L4").Value = 1.45
$ws.Range("M4").Value = 2.37
$ws.Range("N4").Value = 2.32
$ws.Range("O4").Value = 1.47
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.18
$ws.Range("R4").Value = 1.98
$ws.Range("S4").Value = 1.65
$ws.Range("T4").Value = 5.8
$ws.Range("U4").Value = 9
$ws.Range("V4").Value = 9.25
$ws.Range("W4").Value = 20
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 40
$ws.Range("Z4").Value = 6.8
$ws.Range("AA4").Value = 5.9
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 110
$ws.Range("AD4").Value = 900
$ws.Range("AE4").Value = 8.25
$ws.Range("AF4").Value = 17
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 50
$ws.Range("AI4").Value = 37
$ws.Range("AJ4").Value = 50

$wb.Save()
